$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.060.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.179.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.177.43'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.15%  '
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("E12").Value = '  +4.00%  '
$ws.Range("E13").Value = '  +3.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.700.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.051.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.180.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("E19").Value = '  +1.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '508.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.61%  '
$ws.Range("E22").Value = '  +3.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +4.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.15%  '
$ws.Range("E31").Value = '  +4.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.93%  '
$ws.Range("E33").Value = '  +3.54%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '486.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0890'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.121'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.98%  '
$ws.Range("E42").Value = '  +6.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0649'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.890.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.117'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +11.92%  '
$ws.Range("E51").Value = '  +3.72%  '
